$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.222.85"
$ws.Range("D3").Value = "1.906.00"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'307.87"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.5266"
$ws.Range("E7").Value = "  +0.63%  "
$ws.Range("D8").Value = "'0.3821"
$ws.Range("E8").Value = "  +1.62%  "
$ws.Range("D9").Value = "'0.07308"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").Value = "'0.9063"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("D13").Value = "'96.11"
$ws.Range("E13").Value = "  +1.10%  "
$ws.Range("D14").Value = "'5.371"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "1.764.68"
$ws.Range("E15").Value = "  -7.61%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "'0.000008683"
$ws.Range("E17").Value = "  +0.42%  "
$ws.Range("D18").Value = "'14.76"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").Value = "27.262.43"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "'5.122"
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "'10.83"
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").Value = "'6.504"
$ws.Range("E23").Value = "  +1.22%  "
$ws.Range("D24").Value = "'2.346"
$ws.Range("E24").Value = "  +2.87%  "
$ws.Range("D25").Value = "'150.24"
$ws.Range("E25").Value = "  +1.87%  "
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("D27").Value = "'1.742"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'116.73"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'4.856"
$ws.Range("E29").Value = "  +0.89%  "
$ws.Range("D30").Value = "'4.882"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'0.09225"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "'0.8235"
$ws.Range("E32").Value = "  +1.81%  "
$ws.Range("D33").Value = "'0.05082"
$ws.Range("E33").Value = "  +0.49%  "
$ws.Range("D34").Value = "'1.233"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "'2.995"
$ws.Range("E35").Value = "  +1.59%  "
$ws.Range("D36").Value = "'3.363"
$ws.Range("E36").Value = "  -2.13%  "
$ws.Range("D37").Value = "'2.720"
$ws.Range("E37").Value = "  +3.82%  "
$ws.Range("D38").Value = "'0.5742"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").Value = "'0.02002"
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "'1.083"
$ws.Range("E40").Value = "  +0.64%  "
$ws.Range("D41").Value = "'9.043"
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("D42").Value = "'6.611"
$ws.Range("E42").Value = "  -0.30%  "
$ws.Range("D43").Value = "'117.11"
$ws.Range("E43").Value = "  +0.49%  "
$ws.Range("E44").Value = "  +0.75%  "
$ws.Range("D45").Value = "'0.4939"
$ws.Range("E45").Value = "  +1.81%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").Value = "'1.643"
$ws.Range("E48").Value = "  +1.69%  "
$ws.Range("D49").Value = "'38.64"
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("D50").Value = "'64.22"
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +0.35%  "
